# Update handback/handoff timestamps in the Generate Report for Handback step.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on Overview and the matching
# "Correspond Handoff Datetime" on de-de share the same handoff timestamp.
$wsOverview.Range("G2").Value = "2016-08-15 17:02:12"
$wsDeDe.Range("H2").Value = "2016-08-15 17:02:12"

# zh-cn handoff / handback datetimes for the b84f2937 file.
$wsZhCn.Range("H2").Value = "2016-08-15 17:02:02"
$wsZhCn.Range("K2").Value = "2016-08-15 17:02:29"

# de-de handback datetime for the b84f2937 file.
$wsDeDe.Range("K2").Value = "2016-08-15 17:02:37"
